$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2023-12-27 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-28 Thursday", 2) | Out-Null

# Update the answer table (20 rows x 5 columns), addressed positionally
# to correctly handle duplicate cell values that map to different targets.
$table = $d.Tables.Item(1)

$table.Cell(1, 1).Range.Text = "77-72=5"
$table.Cell(1, 2).Range.Text = "16+41=57"
$table.Cell(1, 3).Range.Text = "29+34=63"
$table.Cell(1, 4).Range.Text = "46+23=69"
$table.Cell(1, 5).Range.Text = "82-18=64"

$table.Cell(2, 1).Range.Text = "96-40=56"
$table.Cell(2, 2).Range.Text = "23+40=63"
$table.Cell(2, 3).Range.Text = "55+42=97"
$table.Cell(2, 4).Range.Text = "25+51=76"
$table.Cell(2, 5).Range.Text = "72-70=2"

$table.Cell(3, 1).Range.Text = "45+5=50"
$table.Cell(3, 2).Range.Text = "56+7=63"
$table.Cell(3, 3).Range.Text = "74-36=38"
$table.Cell(3, 4).Range.Text = "39-39=0"
$table.Cell(3, 5).Range.Text = "27+19=46"

$table.Cell(4, 1).Range.Text = "39+18=57"
$table.Cell(4, 2).Range.Text = "64-12=52"
$table.Cell(4, 3).Range.Text = "27+64=91"
$table.Cell(4, 4).Range.Text = "88+7=95"
$table.Cell(4, 5).Range.Text = "40+48=88"

$table.Cell(5, 1).Range.Text = "57+35=92"
$table.Cell(5, 2).Range.Text = "44+50=94"
$table.Cell(5, 3).Range.Text = "91-10=81"
$table.Cell(5, 4).Range.Text = "89-88=1"
$table.Cell(5, 5).Range.Text = "26-26=0"

$table.Cell(6, 1).Range.Text = "60+14=74"
$table.Cell(6, 2).Range.Text = "27-4=23"
$table.Cell(6, 3).Range.Text = "70-22=48"
$table.Cell(6, 4).Range.Text = "88-23=65"
$table.Cell(6, 5).Range.Text = "54+20=74"

$table.Cell(7, 1).Range.Text = "94-17=77"
$table.Cell(7, 2).Range.Text = "63-23=40"
$table.Cell(7, 3).Range.Text = "5+57=62"
$table.Cell(7, 4).Range.Text = "23+19=42"
$table.Cell(7, 5).Range.Text = "3+77=80"

$table.Cell(8, 1).Range.Text = "90+7=97"
$table.Cell(8, 2).Range.Text = "57-20=37"
$table.Cell(8, 3).Range.Text = "89-21=68"
$table.Cell(8, 4).Range.Text = "41-20=21"
$table.Cell(8, 5).Range.Text = "39-25=14"

$table.Cell(9, 1).Range.Text = "32+42=74"
$table.Cell(9, 2).Range.Text = "76+12=88"
$table.Cell(9, 3).Range.Text = "43-35=8"
$table.Cell(9, 4).Range.Text = "15+58=73"
$table.Cell(9, 5).Range.Text = "42+15=57"

$table.Cell(10, 1).Range.Text = "90-48=42"
$table.Cell(10, 2).Range.Text = "43-15=28"
$table.Cell(10, 3).Range.Text = "1+78=79"
$table.Cell(10, 4).Range.Text = "65-13=52"
$table.Cell(10, 5).Range.Text = "82-13=69"

$table.Cell(11, 1).Range.Text = "94-43=51"
$table.Cell(11, 2).Range.Text = "49-32=17"
$table.Cell(11, 3).Range.Text = "31+25=56"
$table.Cell(11, 4).Range.Text = "68+4=72"
$table.Cell(11, 5).Range.Text = "13+79=92"

$table.Cell(12, 1).Range.Text = "29-22=7"
$table.Cell(12, 2).Range.Text = "88+8=96"
$table.Cell(12, 3).Range.Text = "14+13=27"
$table.Cell(12, 4).Range.Text = "54-8=46"
$table.Cell(12, 5).Range.Text = "25-24=1"

$table.Cell(13, 1).Range.Text = "16+76=92"
$table.Cell(13, 2).Range.Text = "17+73=90"
$table.Cell(13, 3).Range.Text = "63+34=97"
$table.Cell(13, 4).Range.Text = "74-0=74"
$table.Cell(13, 5).Range.Text = "52-29=23"

$table.Cell(14, 1).Range.Text = "77+14=91"
$table.Cell(14, 2).Range.Text = "50-8=42"
$table.Cell(14, 3).Range.Text = "37-32=5"
$table.Cell(14, 4).Range.Text = "43+5=48"
$table.Cell(14, 5).Range.Text = "9+79=88"

$table.Cell(15, 1).Range.Text = "97-90=7"
$table.Cell(15, 2).Range.Text = "6+11=17"
$table.Cell(15, 3).Range.Text = "68-10=58"
$table.Cell(15, 4).Range.Text = "88-55=33"
$table.Cell(15, 5).Range.Text = "53-45=8"

$table.Cell(16, 1).Range.Text = "98-7=91"
$table.Cell(16, 2).Range.Text = "80-30=50"
$table.Cell(16, 3).Range.Text = "28-10=18"
$table.Cell(16, 4).Range.Text = "3+96=99"
$table.Cell(16, 5).Range.Text = "0+87=87"

$table.Cell(17, 1).Range.Text = "5+21=26"
$table.Cell(17, 2).Range.Text = "65-65=0"
$table.Cell(17, 3).Range.Text = "57+0=57"
$table.Cell(17, 4).Range.Text = "34+64=98"
$table.Cell(17, 5).Range.Text = "86-82=4"

$table.Cell(18, 1).Range.Text = "66-65=1"
$table.Cell(18, 2).Range.Text = "68+29=97"
$table.Cell(18, 3).Range.Text = "66-8=58"
$table.Cell(18, 4).Range.Text = "25+31=56"
$table.Cell(18, 5).Range.Text = "29+50=79"

$table.Cell(19, 1).Range.Text = "94-36=58"
$table.Cell(19, 2).Range.Text = "15+57=72"
$table.Cell(19, 3).Range.Text = "98-79=19"
$table.Cell(19, 4).Range.Text = "58+26=84"
$table.Cell(19, 5).Range.Text = "95-5=90"

$table.Cell(20, 1).Range.Text = "98-50=48"
$table.Cell(20, 2).Range.Text = "71-27=44"
$table.Cell(20, 3).Range.Text = "19+7=26"
$table.Cell(20, 4).Range.Text = "20+59=79"
$table.Cell(20, 5).Range.Text = "27+52=79"

